$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 31
$ws.Cells.Item($row, 1).Value = 56
$ws.Cells.Item($row, 2).Value = "changes updated"
$ws.Cells.Item($row, 3).Value = "riya-morankar"
$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "edit1 to main"

# Force the date column to be stored as literal text (matching the rest of
# the sheet, which uses plain text dates) instead of being auto-converted
# into an Excel date serial number.
$ws.Cells.Item($row, 6).NumberFormat = "@"
$ws.Cells.Item($row, 6).Value = "2025-06-23"
$ws.Cells.Item($row, 6).Style = "Normal"
